$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D5","D6","D11","D14","D23","D28","D30","D31","D36","D38","D43","D45")
foreach ($ref in $textCells) { $ws.Range($ref).NumberFormat = "@" }

$ws.Range("D2").Value = "67.085.95"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "2.466.82"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "582.07"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").Value = "173.72"
$ws.Range("E6").Value = "  +2.58%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -0.52%  "
$ws.Range("E9").Value = "  +1.33%  "
$ws.Range("E10").Value = "  +0.19%  "
$ws.Range("D11").Value = "4.93"
$ws.Range("E11").Value = "  -0.39%  "
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("D14").Value = "25.33"
$ws.Range("E14").Value = "  -1.01%  "
$ws.Range("D15").Value = "66.689.88"
$ws.Range("E15").Value = "  +0.11%  "
$ws.Range("D17").Value = "2.431.16"
$ws.Range("E17").Value = "  -0.69%  "
$ws.Range("E18").Value = "  -2.29%  "
$ws.Range("E19").Value = "  -1.77%  "
$ws.Range("E20").Value = "  -1.27%  "
$ws.Range("E21").Value = "  -0.35%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").Value = "69.34"
$ws.Range("E23").Value = "  +0.48%  "
$ws.Range("E24").Value = "  -1.44%  "
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("E26").Value = "  -1.26%  "
$ws.Range("D27").Value = "2.596.02"
$ws.Range("E27").Value = "  +0.38%  "
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("D29").Value = "0.0₃0896"
$ws.Range("E29").Value = "  -0.90%  "
$ws.Range("D30").Value = "497.26"
$ws.Range("E30").Value = "  -3.88%  "
$ws.Range("D31").Value = "7.71"
$ws.Range("E31").Value = "  -0.46%  "
$ws.Range("E32").Value = "  -0.67%  "
$ws.Range("E33").Value = "  -1.41%  "
$ws.Range("E35").Value = "  +2.06%  "
$ws.Range("D36").Value = "161.80"
$ws.Range("E36").Value = "  +1.96%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").Value = "18.11"
$ws.Range("E38").Value = "  -1.37%  "
$ws.Range("E39").Value = "  -2.12%  "
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("E41").Value = "  +0.76%  "
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("D43").Value = "4.80"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("D45").Value = "142.23"
$ws.Range("E45").Value = "  +0.80%  "
$ws.Range("E46").Value = "  +0.43%  "
$ws.Range("E47").Value = "  -1.37%  "
$ws.Range("D48").Value = "0.0₆0252"
$ws.Range("E48").Value = "  -0.92%  "
$ws.Range("E49").Value = "  +0.87%  "
$ws.Range("E50").Value = "  -1.78%  "
$ws.Range("E51").Value = "  -0.02%  "

foreach ($ref in $textCells) { $ws.Range($ref).Style = "Normal" }
